# Update the betting-odds values in row 2 to reflect the latest FlashScore
# data refresh (Atualizando o arquivo XLSX).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1.6
    "H2"  = 3.5
    "I2"  = 6.5
    "J2"  = 2.25
    "K2"  = 2.1
    "M2"  = 1.07
    "N2"  = 9
    "O2"  = 1.36
    "P2"  = 3
    "Q2"  = 2.2
    "R2"  = 1.65
    "X2"  = 6.5
    "Z2"  = 11
    "AA2" = 15
    "AC2" = 7.5
    "AD2" = 7
    "AE2" = 21
    "AF2" = 81
    "AG2" = 13
    "AH2" = 29
    "AI2" = 21
    "AJ2" = 67
    "AN2" = 3.4
    "AO2" = 8.5
    "AQ2" = 29
    "AW2" = 7
    "AX2" = 34
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
